$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.555.82'
$ws.Range('E2').Value = '  +0.61%  '

# Row 3
$ws.Range('D3').Value = '2.341.30'
$ws.Range('E3').Value = '  +0.15%  '

# Row 4
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '549.62'
$ws.Range('E5').Value = '  +0.88%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.51'
$ws.Range('E6').Value = '  -0.22%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.06%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.581'
$ws.Range('E8').Value = '  -0.60%  '

# Row 9
$ws.Range('D9').Value = '2.339.64'
$ws.Range('E9').Value = '  +0.19%  '

# Row 10
$ws.Range('E10').Value = '  +1.48%  '

# Row 11
$ws.Range('E11').Value = '  +2.06%  '

# Row 12
$ws.Range('E12').Value = '  -0.50%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.337'
$ws.Range('E13').Value = '  +1.39%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.92'
$ws.Range('E14').Value = '  +0.47%  '

# Row 15
$ws.Range('D15').Value = '2.756.45'
$ws.Range('E15').Value = '  -0.01%  '

# Row 16
$ws.Range('D16').Value = '60.460.55'
$ws.Range('E16').Value = '  +0.54%  '

# Row 17
$ws.Range('E17').Value = '  +1.48%  '

# Row 18
$ws.Range('D18').Value = '2.339.23'
$ws.Range('E18').Value = '  +0.01%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.70'
$ws.Range('E19').Value = '  +0.99%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.12'
$ws.Range('E20').Value = '  -0.83%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '314.49'
$ws.Range('E21').Value = '  +0.30%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.62'
$ws.Range('E22').Value = '  -2.60%  '

# Row 23
$ws.Range('E23').Value = '  +0.18%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.27'
$ws.Range('E24').Value = '  +1.34%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.171'
$ws.Range('E25').Value = '  -0.63%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.31%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.99'
$ws.Range('E27').Value = '  +1.25%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.40'
$ws.Range('E28').Value = '  +3.65%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.27'
$ws.Range('E29').Value = '  +7.91%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.75'
$ws.Range('E30').Value = '  +0.36%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '171.02'
$ws.Range('E31').Value = '  -0.24%  '

# Row 32
$ws.Range('D32').Value = '0.0₃0738'
$ws.Range('E32').Value = '  +1.43%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.09'
$ws.Range('E33').Value = '  +2.51%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.386'
$ws.Range('E34').Value = '  +1.13%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.37'
$ws.Range('E35').Value = '  -1.95%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.11'
$ws.Range('E36').Value = '  +0.47%  '

# Row 37
$ws.Range('E37').Value = '  +0.01%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.08%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.14'
$ws.Range('E39').Value = '  -0.25%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '329.49'
$ws.Range('E40').Value = '  +2.32%  '

# Row 41
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.54'
$ws.Range('E41').Value = '  +0.99%  '

# Row 42
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.05'
$ws.Range('E42').Value = '  -0.10%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '138.13'
$ws.Range('E43').Value = '  -2.15%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.52'
$ws.Range('E44').Value = '  +1.86%  '

# Row 45
$ws.Range('E45').Value = '  +0.66%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.34'
$ws.Range('E46').Value = '  -0.88%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.570'
$ws.Range('E47').Value = '  +1.67%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0500'
$ws.Range('E48').Value = '  +0.51%  '

# Row 49
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0216'
$ws.Range('E49').Value = '  +1.80%  '

# Row 50
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0222'
$ws.Range('E50').Value = '  +7.37%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.94'
$ws.Range('E51').Value = '  -0.84%  '

